# Apply weekly price-update edit to "Fruta, Vega Modelo de Temuco - Níspero"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows that changed (columns D, L, M, N, O, P, Q, R, S, T)
$updates = @(
    @{ Row = 2;  D = 44488; L = "Primera"; M = 100; N = 12000; O = 12000; P = 12000; Q = "$/bandeja 5 kilos";  R = "La Ligua";               S = 2400; T = 5 },
    @{ Row = 4;  D = 44496; L = "Primera"; M = 55;  N = 28000; O = 28000; P = 28000; Q = "$/bandeja 10 kilos"; R = "Provincia de Quillota";   S = 2800; T = 10 },
    @{ Row = 5;  D = 44519; L = "Primera"; M = 30;  N = 28000; O = 28000; P = 28000; Q = "$/bandeja 10 kilos"; R = "Provincia de Quillota";   S = 2800; T = 10 },
    @{ Row = 7;  D = 44166; L = "Segunda"; M = 20;  N = 12000; O = 12000; P = 12000; Q = "$/caja 18 kilos";    R = "La Ligua";               S = 667;  T = 18 },
    @{ Row = 8;  D = 44515; L = "Primera"; M = 80;  N = 28000; O = 28000; P = 28000; Q = "$/bandeja 10 kilos"; R = "Provincia de Los Andes"; S = 2800; T = 10 },
    @{ Row = 9;  D = 44511; L = "Primera"; M = 45;  N = 28000; O = 28000; P = 28000; Q = "$/bandeja 10 kilos"; R = "Provincia de Los Andes"; S = 2800; T = 10 },
    @{ Row = 10; D = 44511; L = "Primera"; M = 45;  N = 3200;  O = 3200;  P = 3200;  Q = "$/bandeja 10 kilos"; R = "Provincia de Quillota";   S = 320;  T = 10 },
    @{ Row = 11; D = 44483; L = "Primera"; M = 35;  N = 10000; O = 10000; P = 10000; Q = "$/bandeja 5 kilos";  R = "Provincia de Quillota";   S = 2000; T = 5 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value  = $u.D   # D: Fecha
    $ws.Cells.Item($r, 12).Value = $u.L   # L: Calidad
    $ws.Cells.Item($r, 13).Value = $u.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $u.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $u.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $u.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $u.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $u.R   # R: Origen
    $ws.Cells.Item($r, 19).Value = $u.S   # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $u.T   # T: Kg / unidad
}
